# Update the list of NEMO variables in the "pre basic identified missing"
# sheet to match the shaconemo repository revision 204.
#
# The "*Lut" variables (tasLut, tslsiLut, hussLut, hflsLut, hfssLut, rsusLut,
# rlusLut, sweLut, fahLut) in rows 144-152 previously all shared the comment
# "Can not be produced by LPJ-GUESS: H-TESSEL?" attributed to "David Warlind".
# They are now updated with a clearer comment and a second comment author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newComment = "Can not be produced by either LPJ-GUESS or H-TESSEL."
$newAuthor  = "David Warlind & Andrea Alessandri"

$firstRow = 144
$lastRow  = 152

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $newComment
    $ws.Cells.Item($r, 7).Value = $newAuthor
}

# Reflect the editor's on-screen state at the time of the commit: scrolled
# down and with the variable column (C) selected over the remaining rows.
[void]$ws.Range("C191:C404").Select()
